{"js": "// Update the date line and the 100 addition/subtraction problems in the\n// table, keeping each cell's existing run/paragraph formatting intact.\n// Each paragraph/cell holds exactly one run of text that must be swapped\n// for a new expression (or, for the first paragraph, a new date string).\n\nconst newDate = \"2024-02-01 Thursday\";\n\n// Row-major replacement values for the 5-column x 20-row table, in the\n// exact order the cells appear in the document (top-left to bottom-right).\nconst newValues = [\n  [\"95+4=\", \"17+53=\", \"89-3=\", \"33+50=\", \"71-36=\"],\n  [\"44-18=\", \"70+4=\", \"88-43=\", \"21-9=\", \"20+18=\"],\n  [\"49-33=\", \"17-4=\", \"87+1=\", \"30+30=\", \"47-43=\"],\n  [\"58+18=\", \"25+19=\", \"52-12=\", \"41-16=\", \"63-55=\"],\n  [\"13+70=\", \"88-22=\", \"22+33=\", \"63+27=\", \"32+7=\"],\n  [\"14+4=\", \"21+55=\", \"88+8=\", \"29+55=\", \"94-77=\"],\n  [\"72+3=\", \"15+63=\", \"23+29=\", \"53-6=\", \"62-29=\"],\n  [\"23+62=\", \"42+1=\", \"75-67=\", \"31+42=\", \"60-20=\"],\n  [\"29-8=\", \"10+57=\", \"56-48=\", \"72-46=\", \"39+44=\"],\n  [\"3+72=\", \"49-13=\", \"72-6=\", \"87-87=\", \"30-18=\"],\n  [\"30+51=\", \"59+9=\", \"3+35=\", \"18+5=\", \"37-13=\"],\n  [\"67-44=\", \"73-36=\", \"77-45=\", \"79-11=\", \"44-24=\"],\n  [\"72+3=\", \"54-22=\", \"53+20=\", \"93-22=\", \"7-2=\"],\n  [\"31+44=\", \"36-9=\", \"31-30=\", \"12+51=\", \"87+5=\"],\n  [\"21+49=\", \"74-51=\", \"54+36=\", \"62-51=\", \"77-25=\"],\n  [\"36-36=\", \"94-4=\", \"38-12=\", \"90-19=\", \"0+93=\"],\n  [\"44-0=\", \"82+17=\", \"26-5=\", \"95-74=\", \"70-7=\"],\n  [\"25+64=\", \"1+57=\", \"23+34=\", \"70-69=\", \"99-59=\"],\n  [\"81+4=\", \"53+32=\", \"64-28=\", \"38-12=\", \"12+6=\"],\n  [\"25+74=\", \"89-73=\", \"79-61=\", \"18-8=\", \"79-26=\"],\n];\n\n// 1) Update the date paragraph (first paragraph in the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.getRange().insertText(newDate, Word.InsertLocation.replace);\n\n// 2) Update every cell of the (single) table with its new expression.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (let r = 0; r < newValues.length; r++) {\n  const rowValues = newValues[r];\n  for (let c = 0; c < rowValues.length; c++) {\n    const cell = table.getCell(r, c);\n    const para = cell.body.paragraphs.getFirst();\n    para.getRange().insertText(rowValues[c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the date line (first paragraph in the body).\n$d.Paragraphs.Item(1).Range.Text = \"2024-02-01 Thursday\"\n\n# 2) Update every cell of the (single) table with its new expression,\n#    preserving each cell's existing paragraph/run formatting.\n$t = $d.Tables.Item(1)\n\n$t.Cell(1,1).Range.Text = \"95+4=\"\n$t.Cell(1,2).Range.Text = \"17+53=\"\n$t.Cell(1,3).Range.Text = \"89-3=\"\n$t.Cell(1,4).Range.Text = \"33+50=\"\n$t.Cell(1,5).Range.Text = \"71-36=\"\n\n$t.Cell(2,1).Range.Text = \"44-18=\"\n$t.Cell(2,2).Range.Text = \"70+4=\"\n$t.Cell(2,3).Range.Text = \"88-43=\"\n$t.Cell(2,4).Range.Text = \"21-9=\"\n$t.Cell(2,5).Range.Text = \"20+18=\"\n\n$t.Cell(3,1).Range.Text = \"49-33=\"\n$t.Cell(3,2).Range.Text = \"17-4=\"\n$t.Cell(3,3).Range.Text = \"87+1=\"\n$t.Cell(3,4).Range.Text = \"30+30=\"\n$t.Cell(3,5).Range.Text = \"47-43=\"\n\n$t.Cell(4,1).Range.Text = \"58+18=\"\n$t.Cell(4,2).Range.Text = \"25+19=\"\n$t.Cell(4,3).Range.Text = \"52-12=\"\n$t.Cell(4,4).Range.Text = \"41-16=\"\n$t.Cell(4,5).Range.Text = \"63-55=\"\n\n$t.Cell(5,1).Range.Text = \"13+70=\"\n$t.Cell(5,2).Range.Text = \"88-22=\"\n$t.Cell(5,3).Range.Text = \"22+33=\"\n$t.Cell(5,4).Range.Text = \"63+27=\"\n$t.Cell(5,5).Range.Text = \"32+7=\"\n\n$t.Cell(6,1).Range.Text = \"14+4=\"\n$t.Cell(6,2).Range.Text = \"21+55=\"\n$t.Cell(6,3).Range.Text = \"88+8=\"\n$t.Cell(6,4).Range.Text = \"29+55=\"\n$t.Cell(6,5).Range.Text = \"94-77=\"\n\n$t.Cell(7,1).Range.Text = \"72+3=\"\n$t.Cell(7,2).Range.Text = \"15+63=\"\n$t.Cell(7,3).Range.Text = \"23+29=\"\n$t.Cell(7,4).Range.Text = \"53-6=\"\n$t.Cell(7,5).Range.Text = \"62-29=\"\n\n$t.Cell(8,1).Range.Text = \"23+62=\"\n$t.Cell(8,2).Range.Text = \"42+1=\"\n$t.Cell(8,3).Range.Text = \"75-67=\"\n$t.Cell(8,4).Range.Text = \"31+42=\"\n$t.Cell(8,5).Range.Text = \"60-20=\"\n\n$t.Cell(9,1).Range.Text = \"29-8=\"\n$t.Cell(9,2).Range.Text = \"10+57=\"\n$t.Cell(9,3).Range.Text = \"56-48=\"\n$t.Cell(9,4).Range.Text = \"72-46=\"\n$t.Cell(9,5).Range.Text = \"39+44=\"\n\n$t.Cell(10,1).Range.Text = \"3+72=\"\n$t.Cell(10,2).Range.Text = \"49-13=\"\n$t.Cell(10,3).Range.Text = \"72-6=\"\n$t.Cell(10,4).Range.Text = \"87-87=\"\n$t.Cell(10,5).Range.Text = \"30-18=\"\n\n$t.Cell(11,1).Range.Text = \"30+51=\"\n$t.Cell(11,2).Range.Text = \"59+9=\"\n$t.Cell(11,3).Range.Text = \"3+35=\"\n$t.Cell(11,4).Range.Text = \"18+5=\"\n$t.Cell(11,5).Range.Text = \"37-13=\"\n\n$t.Cell(12,1).Range.Text = \"67-44=\"\n$t.Cell(12,2).Range.Text = \"73-36=\"\n$t.Cell(12,3).Range.Text = \"77-45=\"\n$t.Cell(12,4).Range.Text = \"79-11=\"\n$t.Cell(12,5).Range.Text = \"44-24=\"\n\n$t.Cell(13,1).Range.Text = \"72+3=\"\n$t.Cell(13,2).Range.Text = \"54-22=\"\n$t.Cell(13,3).Range.Text = \"53+20=\"\n$t.Cell(13,4).Range.Text = \"93-22=\"\n$t.Cell(13,5).Range.Text = \"7-2=\"\n\n$t.Cell(14,1).Range.Text = \"31+44=\"\n$t.Cell(14,2).Range.Text = \"36-9=\"\n$t.Cell(14,3).Range.Text = \"31-30=\"\n$t.Cell(14,4).Range.Text = \"12+51=\"\n$t.Cell(14,5).Range.Text = \"87+5=\"\n\n$t.Cell(15,1).Range.Text = \"21+49=\"\n$t.Cell(15,2).Range.Text = \"74-51=\"\n$t.Cell(15,3).Range.Text = \"54+36=\"\n$t.Cell(15,4).Range.Text = \"62-51=\"\n$t.Cell(15,5).Range.Text = \"77-25=\"\n\n$t.Cell(16,1).Range.Text = \"36-36=\"\n$t.Cell(16,2).Range.Text = \"94-4=\"\n$t.Cell(16,3).Range.Text = \"38-12=\"\n$t.Cell(16,4).Range.Text = \"90-19=\"\n$t.Cell(16,5).Range.Text = \"0+93=\"\n\n$t.Cell(17,1).Range.Text = \"44-0=\"\n$t.Cell(17,2).Range.Text = \"82+17=\"\n$t.Cell(17,3).Range.Text = \"26-5=\"\n$t.Cell(17,4).Range.Text = \"95-74=\"\n$t.Cell(17,5).Range.Text = \"70-7=\"\n\n$t.Cell(18,1).Range.Text = \"25+64=\"\n$t.Cell(18,2).Range.Text = \"1+57=\"\n$t.Cell(18,3).Range.Text = \"23+34=\"\n$t.Cell(18,4).Range.Text = \"70-69=\"\n$t.Cell(18,5).Range.Text = \"99-59=\"\n\n$t.Cell(19,1).Range.Text = \"81+4=\"\n$t.Cell(19,2).Range.Text = \"53+32=\"\n$t.Cell(19,3).Range.Text = \"64-28=\"\n$t.Cell(19,4).Range.Text = \"38-12=\"\n$t.Cell(19,5).Range.Text = \"12+6=\"\n\n$t.Cell(20,1).Range.Text = \"25+74=\"\n$t.Cell(20,2).Range.Text = \"89-73=\"\n$t.Cell(20,3).Range.Text = \"79-61=\"\n$t.Cell(20,4).Range.Text = \"18-8=\"\n$t.Cell(20,5).Range.Text = \"79-26=\"\n\n"}
